$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3031945824623108
$ws.Range("B1").Value = 0.8837385177612305
$ws.Range("C1").Value = 2.845041513442993
$ws.Range("D1").Value = 2.687310934066772
$ws.Range("E1").Value = 1.115850806236267
